# Add a new "Save" column (column H) to the s_vals worksheet, mirroring the
# header style used by the existing summary column (G / "sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy formatting from the existing "sum" header (G1) so the new
# "Save" header (H1) matches the bold/bordered/centered style used by the
# other headers, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column, one flag (0/1) per row.
$saveValues = @(0, 0, 1, 1, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = $false
